$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: text cells in this sheet are formatted with a "quote-prefix" cell
# style (they show up as style index 3 with t="s" in the OOXML). Assigning
# a plain string via .Value flips the cell to the non-quote-prefixed
# sibling style, so a leading "'" is used to force text entry and keep the
# original quote-prefixed styling intact (the leading apostrophe itself is
# never stored in the cell's value).

# --- Row 2 (0402 1uF capacitor) : supplier stock counts drop 4428 -> 4328 ---
$ws.Range("O2").Value = 4328
$ws.Range("R2").Value = 4328

# --- Row 3 (0805 10uF capacitor) : supplier stock counts drop 2599 -> 2593 ---
$ws.Range("O3").Value = 2593
$ws.Range("R3").Value = 2593

# --- Row 4 (0402 10pF capacitor) : supplier stock counts drop 28941 -> 28567 ---
$ws.Range("O4").Value = 28567
$ws.Range("R4").Value = 28567

# --- Row 5 (0402 100nF capacitor) : supplier stock counts drop 10764 -> 10564 ---
$ws.Range("O5").Value = 10564
$ws.Range("R5").Value = 10564

# --- Row 6 (1N4148WT diode) : swap the two Digi-Key order codes, restock ---
$ws.Range("L6").Value = "'1N4148WTDKR-ND"
$ws.Range("Q6").Value = "'1N4148WTCT-ND"
$ws.Range("O6").Value = 305174
$ws.Range("R6").Value = 305174

# --- Row 7 (TVS diode D72) : fill in the missing description + Line# marker,
#     restock ---
$ws.Range("A7").Value = "'"
$ws.Range("E7").Value = "'WE-TVS TVS Diode, High Speed Series, SOT23-6L, VDC = 5V"
$ws.Rows.Item(7).RowHeight = 57.6
$ws.Range("O7").Value = 1155
$ws.Range("R7").Value = 1155

# --- Row 8 (USB connector) : supplier stock counts drop 2611 -> 2561 ---
$ws.Range("O8").Value = 2561
$ws.Range("R8").Value = 2561

# --- Row 10 (22R resistor) : supplier stock counts drop 37993 -> 28260 ---
$ws.Range("O10").Value = 28260
$ws.Range("R10").Value = 28260

# --- Row 11 (4.7K resistor) : supplier stock counts rise 618223 -> 717639 ---
$ws.Range("O11").Value = 717639
$ws.Range("R11").Value = 717639

# --- Row 12 (tactile switch) : supplier stock counts drop 17053 -> 15851 ---
$ws.Range("O12").Value = 15851
$ws.Range("R12").Value = 15851

# --- Row 13 (ATmega32U4-MU) : re-sourced from Digi-Key to Farnell / RS ---
$ws.Range("H13").Value = "'Microchip"
$ws.Range("J13").Value = "'Unknown"
$ws.Range("K13").Value = "'Farnell"
$ws.Range("L13").Value = "'2425127"
$ws.Range("M13").Value = 4.45
$ws.Range("N13").Value = 89.01
$ws.Range("O13").Value = 5516
$ws.Range("P13").Value = "'RSComponents"
$ws.Range("Q13").Value = "'1310290"
$ws.Range("R13").Value = 179
$ws.Range("S13").Value = 78.24

# --- Row 14 (16MHz crystal) : supplier stock counts rise 148 -> 838 ---
$ws.Range("O14").Value = 838
$ws.Range("R14").Value = 838
